$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D:G range to text format to preserve exact string formatting (percent/decimal
# look-alikes) as literal text, matching the inlineStr source values.
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "301.73"
$ws.Range("E2").Value = "1.96%"
$ws.Range("G2").Value = "17"
$ws.Range("D3").Value = "44.14"
$ws.Range("E3").Value = "6.78%"
$ws.Range("G3").Value = "17"
$ws.Range("D4").Value = "5.083"
$ws.Range("E4").Value = "1.24%"
$ws.Range("G4").Value = "17"
$ws.Range("D5").Value = "0.07701"
$ws.Range("E5").Value = "3.26%"
$ws.Range("G5").Value = "17"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "4.425"
$ws.Range("E6").Value = "1.55%"
$ws.Range("G6").Value = "17"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.616"
$ws.Range("E7").Value = "2.50%"
$ws.Range("G7").Value = "17"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "1.036"
$ws.Range("E8").Value = "11.54%"
$ws.Range("G8").Value = "17"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1275"
$ws.Range("E9").Value = "8.15%"
$ws.Range("G9").Value = "17"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1867"
$ws.Range("E10").Value = "3.20%"
$ws.Range("G10").Value = "17"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.09266"
$ws.Range("E11").Value = "5.09%"
$ws.Range("G11").Value = "17"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.04175"
$ws.Range("E12").Value = "-0.31%"
$ws.Range("G12").Value = "17"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.1047"
$ws.Range("E13").Value = "-0.23%"
$ws.Range("G13").Value = "17"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001281"
$ws.Range("E14").Value = "0.62%"
$ws.Range("G14").Value = "17"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005772"
$ws.Range("E15").Value = "-1.99%"
$ws.Range("G15").Value = "17"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "0.007430"
$ws.Range("E16").Value = "1,897.29%"
$ws.Range("G16").Value = "17"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.346"
$ws.Range("E17").Value = "-0.25%"
$ws.Range("G17").Value = "17"
$ws.Range("D18").Value = "2.331"
$ws.Range("E18").Value = "-3.38%"
$ws.Range("G18").Value = "17"
$ws.Range("E19").Value = "1.92%"
$ws.Range("G19").Value = "17"
$ws.Range("D20").Value = "8.638"
$ws.Range("E20").Value = "9.44%"
$ws.Range("G20").Value = "17"
$ws.Range("D21").Value = "0.1400"
$ws.Range("E21").Value = "-0.64%"
$ws.Range("G21").Value = "17"
$ws.Range("D22").Value = "0.3178"
$ws.Range("E22").Value = "6.94%"
$ws.Range("G22").Value = "17"
$ws.Range("D23").Value = "0.04200"
$ws.Range("E23").Value = "4.32%"
$ws.Range("G23").Value = "17"
$ws.Range("D24").Value = "0.001286"
$ws.Range("E24").Value = "1.63%"
$ws.Range("G24").Value = "17"
$ws.Range("D25").Value = "0.004470"
$ws.Range("E25").Value = "15.76%"
$ws.Range("G25").Value = "17"
$ws.Range("D26").Value = "0.0001350"
$ws.Range("E26").Value = "9.78%"
$ws.Range("G26").Value = "17"
$ws.Range("G27").Value = "17"
$ws.Range("G28").Value = "17"
$ws.Range("G29").Value = "17"
$ws.Range("G30").Value = "17"
$ws.Range("G31").Value = "17"
$ws.Range("G32").Value = "17"
$ws.Range("G33").Value = "17"
$ws.Range("G34").Value = "17"
$ws.Range("G35").Value = "17"
$ws.Range("G36").Value = "17"
$ws.Range("G37").Value = "17"
$ws.Range("D38").Value = "0.02487"
$ws.Range("E38").Value = "4.08%"
$ws.Range("G38").Value = "17"
$ws.Range("D39").Value = "0.05298"
$ws.Range("E39").Value = "1.89%"
$ws.Range("G39").Value = "17"
$ws.Range("D40").Value = "0.005979"
$ws.Range("E40").Value = "-8.95%"
$ws.Range("G40").Value = "17"
$ws.Range("D41").Value = "0.007717"
$ws.Range("E41").Value = "-0.88%"
$ws.Range("G41").Value = "17"
$ws.Range("E42").Value = "2.63%"
$ws.Range("G42").Value = "17"
$ws.Range("D43").Value = "0.007352"
$ws.Range("E43").Value = "-0.40%"
$ws.Range("G43").Value = "17"
$ws.Range("D44").Value = "0.007568"
$ws.Range("E44").Value = "-2.97%"
$ws.Range("G44").Value = "17"
$ws.Range("D45").Value = "0.2999"
$ws.Range("E45").Value = "-6.73%"
$ws.Range("G45").Value = "17"
$ws.Range("D46").Value = "0.00006664"
$ws.Range("E46").Value = "6.75%"
$ws.Range("G46").Value = "17"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "0.03%"
$ws.Range("G47").Value = "17"
$ws.Range("D48").Value = "0.04188"
$ws.Range("E48").Value = "-9.19%"
$ws.Range("G48").Value = "17"
$ws.Range("E49").Value = "0.06%"
$ws.Range("G49").Value = "17"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "0.03%"
$ws.Range("G50").Value = "17"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "0.03%"
$ws.Range("G51").Value = "17"

# Restore default style (style index 0) on the touched range so no stray
# cell-level number-format styling is left behind.
$ws.Range("D2:G51").Style = "Normal"

